$p = $ppt.ActivePresentation

# Locate the slide with SlideID 261 (7th slide in the deck)
$s = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    if ($candidate.SlideID -eq 261) {
        $s = $candidate
        break
    }
}
if ($s -eq $null) {
    $s = $p.Slides.Item(7)
}

$shp = $s.Shapes.Item("TextBox 22")

# Reposition / resize the "Output Probabilities" textbox
# (EMU -> points, nudged slightly so the float32 COM properties
#  round-trip back to the exact target EMU values)
$shp.Left = 654.2912998425197
$shp.Top = 85.91818897637795
$shp.Width = 242.56308086614175
$shp.Height = 24.234419448818898

$tf = $shp.TextFrame
$tr = $tf.TextRange

# Merge the two paragraphs ("Output" / "Probabilities") into a single
# paragraph split across three runs: "output", " ", "probabilities"
$tr.Text = "output probabilities"

$r1 = $tr.Characters(1, 6)
$r1.Font.Name = "+mj-lt"

$r2 = $tr.Characters(7, 1)
$r2.Font.Name = "+mj-lt"

$r3 = $tr.Characters(8, 13)
$r3.Font.Name = "+mj-lt"
